$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2818
$ws.Range("F3").Value = 1139
$ws.Range("F4").Value = 20595
$ws.Range("F5").Value = 96
$ws.Range("F6").Value = 2625
$ws.Range("F7").Value = 784
$ws.Range("F9").Value = 492
$ws.Range("F10").Value = 736
$ws.Range("F11").Value = 272
$ws.Range("F14").Value = 102
$ws.Range("F15").Value = 500
$ws.Range("F17").Value = 244
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 402
$ws.Range("F20").Value = 5
$ws.Range("F22").Value = 25

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 16
$ws.Range("F3").Value = 25
$ws.Range("F5").Value = 317
$ws.Range("F8").Value = 17
$ws.Range("F14").Value = 129

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6089
$ws.Range("F3").Value = 683
$ws.Range("F4").Value = 657
$ws.Range("F5").Value = 1455

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 6089
$ws.Range("F3").Value = 683
$ws.Range("F4").Value = 657
$ws.Range("F5").Value = 1455
$ws.Range("F6").Value = 2818
$ws.Range("F7").Value = 1139
$ws.Range("F8").Value = 20595
$ws.Range("F9").Value = 16
$ws.Range("F10").Value = 25
$ws.Range("F11").Value = 96
$ws.Range("F13").Value = 317
$ws.Range("F14").Value = 2625
$ws.Range("F15").Value = 784
$ws.Range("F19").Value = 492
$ws.Range("F20").Value = 736
$ws.Range("F21").Value = 272
$ws.Range("F25").Value = 17
$ws.Range("F27").Value = 102
$ws.Range("F30").Value = 500
$ws.Range("F34").Value = 244
$ws.Range("F35").Value = 129
$ws.Range("F36").Value = 129
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 402
$ws.Range("F40").Value = 5
$ws.Range("F44").Value = 25

